# "local run for metrics" - bump the n_iter value (column N, row 2) from
# 100 to 1000, and nudge the column C / L widths to match the narrower
# layout the author ended up with after resizing column C (W) by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Core data edit: n_iter 100 -> 1000
$ws.Range("N2").Value = 1000

# Column width tweaks (author narrowed column C / "W", with column L
# picking up a hairline-width rounding nudge as a side effect).
$ws.Columns.Item(3).ColumnWidth = 2.25
$ws.Columns.Item(12).ColumnWidth = 4.5
